$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Categories")

# Column C ("isMissing") rows 2-32 currently hold the text "false" in a
# Text-formatted ("@") cell. Switch those cells to the General number
# format so a numeric value can be stored, then write 0 (isMissing = false).
$dataRange = $ws.Range("C2:C32")
$dataRange.NumberFormat = "General"

for ($r = 2; $r -le 32; $r++) {
    $ws.Cells.Item($r, 3).Value = 0
}

# Move the active selection to C2 (matches the post-edit selection in the file).
$ws.Range("C2").Select()
